$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: change G4 from 15 to 1 (drives the recalculated cascade) ---
$ws1.Range("G4").Value = 1

# --- Add the new second sheet "Foglio2" at the end ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws2.Name = "Foglio2"

# Row 1 - headers
$ws2.Range("A1").Value = "Pietra"
$ws2.Range("C1").Value = "ManaTot"
$ws2.Range("E1").Value = "Cristallo"
$ws2.Range("G1").Value = "ManaTot"

# Bold the two row-1 titles. Build a throwaway named style carrying only the
# "bold font" attribute (no border/alignment) and apply + drop it again so the
# workbook ends up with a plain "applyFont only" cellXf instead of reusing one
# of the pre-existing bordered/centered bold xfs.
$boldStyle = $wb.Styles.Add("__TmpBold")
$boldStyle.Font.Bold = $true
$ws2.Range("A1").Style = "__TmpBold"
$ws2.Range("E1").Style = "__TmpBold"
$wb.Styles("__TmpBold").Delete()

# Row 2
$ws2.Range("A2").Value = "Livello"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Formula = "=SUM(B3:B6)*B2"
$ws2.Range("D2").Value = "MediePerc"
$ws2.Range("E2").Value = "Livello"
$ws2.Range("F2").Value = 5
$ws2.Range("G2").Formula = "=SUM(F3:F7)*F2"
$ws2.Range("K2").Formula = "=1135/5"

# Row 3
$ws2.Range("A3").Value = "Spirito"
$ws2.Range("B3").Value = 20
$ws2.Range("C3").Formula = "=B3*B`$2/C`$2"
$ws2.Range("D3").Formula = "=(C3+G3)/2"
$ws2.Range("E3").Value = "spirito"
$ws2.Range("F3").Value = 10
$ws2.Range("G3").Formula = "=F3*F`$2/G`$2"

# Row 4
$ws2.Range("A4").Value = "Rarita"
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Formula = "=B4*B`$2/C`$2"
$ws2.Range("D4").Formula = "=(C4+G4)/2"
$ws2.Range("E4").Value = "rarita"
$ws2.Range("F4").Value = 0
$ws2.Range("G4").Formula = "=F4*F`$2/G`$2"

# Row 5
$ws2.Range("A5").Value = "attacco"
$ws2.Range("B5").Value = 15
$ws2.Range("C5").Formula = "=B5*B`$2/C`$2"
$ws2.Range("D5").Formula = "=(C5+G5)/2"
$ws2.Range("E5").Value = "attacco"
$ws2.Range("F5").Value = 40
$ws2.Range("G5").Formula = "=F5*F`$2/G`$2"

# Row 6
$ws2.Range("A6").Value = "durezza"
$ws2.Range("B6").Value = 50
$ws2.Range("C6").Formula = "=B6*B`$2/C`$2"
$ws2.Range("D6").Formula = "=(C6+G6)/2"
$ws2.Range("E6").Value = "durezza"
$ws2.Range("F6").Value = 2
$ws2.Range("G6").Formula = "=F6*F`$2/G`$2"

# Row 7
$ws2.Range("E7").Value = "magia"
$ws2.Range("F7").Value = 90
$ws2.Range("G7").Formula = "=F7*F`$2/G`$2"

# Row 9
$ws2.Range("A9").Value = "NuovaPietra"
$ws2.Range("C9").Value = "ManaTot"
$ws2.Range("D9").Value = "Normalizzato"
$ws2.Range("E9").Value = "ManaTot"

# Row 10
$ws2.Range("A10").Value = "Livello"
$ws2.Range("B10").Formula = "=B2"
$ws2.Range("C10").Formula = "=SUM(B11:B14)"
$ws2.Range("D10").Formula = "=150*B10/C10"
$ws2.Range("E10").Formula = "=SUM(D11:D14)"

# Row 11
$ws2.Range("A11").Value = "Spirito"
$ws2.Range("B11").Formula = "=B3*B`$2+D3*G`$2"
$ws2.Range("D11").Formula = "=D`$10*B11"

# Row 12
$ws2.Range("A12").Value = "Rarita"
$ws2.Range("B12").Formula = "=B4*B`$2+C4*G`$2"
$ws2.Range("D12").Formula = "=D`$10*B12"

# Row 13
$ws2.Range("A13").Value = "attacco"
$ws2.Range("B13").Formula = "=B5*B`$2+D5*G`$2"
$ws2.Range("D13").Formula = "=D`$10*B13"

# Row 14
$ws2.Range("A14").Value = "durezza"
$ws2.Range("B14").Formula = "=B6*B`$2+D6*G`$2"
$ws2.Range("D14").Formula = "=D`$10*B14"

# Column widths (match the target workbook as closely as this runtime's
# pixel-quantised ColumnWidth setter allows)
$ws2.Range("C:D").ColumnWidth = 11.5
$ws2.Range("G:G").ColumnWidth = 11

# --- Selections / active sheet ---
$ws1.Range("J12").Select()
$ws2.Range("K3").Select()

Write-Output "done"
